$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, shifting existing rows 202.. down by one.
$ws.Rows(202).Insert()

# Populate the newly inserted row 202 with the new data record.
$ws.Cells.Item(202, 1).Value2  = 4
$ws.Cells.Item(202, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(202, 3).Value2  = "Los Lagos"
$ws.Cells.Item(202, 4).Value2  = 44736
$ws.Cells.Item(202, 5).Value2  = 10
$ws.Cells.Item(202, 6).Value2  = "Fruta"
$ws.Cells.Item(202, 7).Value2  = 100108
$ws.Cells.Item(202, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(202, 9).Value2  = 100108005
$ws.Cells.Item(202, 10).Value2 = "Piña"
$ws.Cells.Item(202, 11).Value2 = "Caramelo"
$ws.Cells.Item(202, 12).Value2 = "Tercera"
$ws.Cells.Item(202, 13).Value2 = 270
$ws.Cells.Item(202, 14).Value2 = 19000
$ws.Cells.Item(202, 15).Value2 = 20000
$ws.Cells.Item(202, 16).Value2 = 19444
$ws.Cells.Item(202, 17).Value2 = "$/caja 16 unidades"
$ws.Cells.Item(202, 18).Value2 = "Ecuador"
$ws.Cells.Item(202, 19).Value2 = 1215
$ws.Cells.Item(202, 20).Value2 = 16
